$wb = $excel.ActiveWorkbook

# This script applies updated market-price derived values (columns H-N)
# across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets, matching a
# scheduled data-refresh run. Values only; no structural changes.

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 96
$ws.Cells.Item(2, 9).Value = 96
$ws.Cells.Item(2, 11).Value = 96
$ws.Cells.Item(2, 13).Value = 17
$ws.Cells.Item(12, 8).Value = 418.4375
$ws.Cells.Item(12, 9).Value = 451.23077
$ws.Cells.Item(12, 11).Value = 451.23077
$ws.Cells.Item(12, 13).Value = -281.23077
$ws.Cells.Item(33, 8).Value = 386.6154
$ws.Cells.Item(33, 9).Value = 363.8
$ws.Cells.Item(33, 11).Value = 363.8
$ws.Cells.Item(33, 13).Value = -134.8
$ws.Cells.Item(41, 8).Value = 531.2
$ws.Cells.Item(41, 10).Value = 346.33334
$ws.Cells.Item(41, 12).Value = 346.33334
$ws.Cells.Item(41, 14).Value = -1226.33334
$ws.Cells.Item(53, 8).Value = 452.33334
$ws.Cells.Item(53, 9).Value = 342.8
$ws.Cells.Item(53, 11).Value = 342.8
$ws.Cells.Item(53, 13).Value = 294.2
$ws.Cells.Item(64, 8).Value = 5251.25
$ws.Cells.Item(64, 10).Value = 6503
$ws.Cells.Item(64, 12).Value = 6503
$ws.Cells.Item(64, 14).Value = -6999
$ws.Cells.Item(67, 8).Value = 5251.25
$ws.Cells.Item(67, 10).Value = 6503
$ws.Cells.Item(67, 12).Value = 6503
$ws.Cells.Item(67, 14).Value = -8219
$ws.Cells.Item(74, 8).Value = 4243.8
$ws.Cells.Item(74, 9).Value = 3493.111
$ws.Cells.Item(74, 11).Value = 3493.111
$ws.Cells.Item(74, 13).Value = -2557.111
$ws.Cells.Item(77, 8).Value = 4243.8
$ws.Cells.Item(77, 9).Value = 3493.111
$ws.Cells.Item(77, 11).Value = 17465.555
$ws.Cells.Item(77, 13).Value = -12785.555
$ws.Cells.Item(111, 8).Value = 2498.5
$ws.Cells.Item(111, 9).Value = 2498.5
$ws.Cells.Item(111, 11).Value = 7495.5
$ws.Cells.Item(111, 13).Value = -4428.5
$ws.Cells.Item(112, 8).Value = 3287.4546
$ws.Cells.Item(112, 10).Value = 3285.111
$ws.Cells.Item(112, 12).Value = 9855.332999999999
$ws.Cells.Item(112, 14).Value = -12071.333
$ws.Cells.Item(116, 8).Value = 9284.833000000001
$ws.Cells.Item(116, 9).Value = 5374.5
$ws.Cells.Item(116, 10).Value = 17105.5
$ws.Cells.Item(116, 11).Value = 5374.5
$ws.Cells.Item(116, 12).Value = 17105.5
$ws.Cells.Item(116, 13).Value = -1932.5
$ws.Cells.Item(116, 14).Value = -23989.5
$ws.Cells.Item(132, 8).Value = 3987.4707
$ws.Cells.Item(132, 9).Value = 3924.1875
$ws.Cells.Item(132, 11).Value = 11772.5625
$ws.Cells.Item(132, 13).Value = -9242.5625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 503.1111
$ws.Cells.Item(4, 9).Value = 504.57144
$ws.Cells.Item(4, 11).Value = 504.57144
$ws.Cells.Item(4, 13).Value = -388.57144
$ws.Cells.Item(63, 8).Value = 3792.8
$ws.Cells.Item(63, 9).Value = 2432.7354
$ws.Cells.Item(63, 10).Value = 11499.833
$ws.Cells.Item(63, 11).Value = 2432.7354
$ws.Cells.Item(63, 12).Value = 11499.833
$ws.Cells.Item(63, 13).Value = -1746.7354
$ws.Cells.Item(63, 14).Value = -12871.833
$ws.Cells.Item(66, 8).Value = 3792.8
$ws.Cells.Item(66, 9).Value = 2432.7354
$ws.Cells.Item(66, 10).Value = 11499.833
$ws.Cells.Item(66, 11).Value = 12163.677
$ws.Cells.Item(66, 12).Value = 57499.165
$ws.Cells.Item(66, 13).Value = -8731.677
$ws.Cells.Item(66, 14).Value = -64363.165
$ws.Cells.Item(94, 8).Value = 40659.57
$ws.Cells.Item(94, 10).Value = 40659.57
$ws.Cells.Item(94, 12).Value = 40659.57
$ws.Cells.Item(94, 14).Value = -42461.57
$ws.Cells.Item(132, 8).Value = 1965.8889
$ws.Cells.Item(132, 10).Value = 2374.25
$ws.Cells.Item(132, 12).Value = 7122.75
$ws.Cells.Item(132, 14).Value = -12182.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 3596.2354
$ws.Cells.Item(86, 9).Value = 1324.9286
$ws.Cells.Item(86, 11).Value = 1324.9286
$ws.Cells.Item(86, 13).Value = -201.9286
$ws.Cells.Item(89, 8).Value = 3596.2354
$ws.Cells.Item(89, 9).Value = 1324.9286
$ws.Cells.Item(89, 11).Value = 6624.643
$ws.Cells.Item(89, 13).Value = -1008.643

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 6361.037
$ws.Cells.Item(22, 9).Value = 1515.5834
$ws.Cells.Item(22, 10).Value = 10237.4
$ws.Cells.Item(22, 11).Value = 1515.5834
$ws.Cells.Item(22, 12).Value = 10237.4
$ws.Cells.Item(22, 13).Value = -1165.5834
$ws.Cells.Item(22, 14).Value = -10937.4
$ws.Cells.Item(86, 8).Value = 4999.5
$ws.Cells.Item(86, 9).Value = 4999.5
$ws.Cells.Item(86, 11).Value = 4999.5
$ws.Cells.Item(86, 13).Value = -3876.5
$ws.Cells.Item(89, 8).Value = 4999.5
$ws.Cells.Item(89, 9).Value = 4999.5
$ws.Cells.Item(89, 11).Value = 24997.5
$ws.Cells.Item(89, 13).Value = -19381.5
$ws.Cells.Item(99, 8).Value = 3999.5
$ws.Cells.Item(99, 9).Value = 3999.5
$ws.Cells.Item(99, 11).Value = 3999.5
$ws.Cells.Item(99, 13).Value = -2501.5
$ws.Cells.Item(126, 8).Value = 3999.5
$ws.Cells.Item(126, 9).Value = 3999.5
$ws.Cells.Item(126, 11).Value = 11998.5
$ws.Cells.Item(126, 13).Value = -9528.5
$ws.Cells.Item(134, 8).Value = 4727.9
$ws.Cells.Item(134, 9).Value = 3871.3333
$ws.Cells.Item(134, 10).Value = 6012.75
$ws.Cells.Item(134, 11).Value = 11613.9999
$ws.Cells.Item(134, 12).Value = 18038.25
$ws.Cells.Item(134, 13).Value = -9078.999899999999
$ws.Cells.Item(134, 14).Value = -23108.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 57.909092
$ws.Cells.Item(12, 9).Value = 64
$ws.Cells.Item(12, 11).Value = 192
$ws.Cells.Item(12, 13).Value = -19
$ws.Cells.Item(14, 8).Value = 1842.5
$ws.Cells.Item(14, 9).Value = 1842.5
$ws.Cells.Item(14, 11).Value = 5527.5
$ws.Cells.Item(14, 13).Value = -5354.5
$ws.Cells.Item(21, 8).Value = 799
$ws.Cells.Item(21, 9).Value = 799
$ws.Cells.Item(21, 11).Value = 2397
$ws.Cells.Item(21, 13).Value = -2224
$ws.Cells.Item(121, 8).Value = 828.1818
$ws.Cells.Item(121, 9).Value = 836.8333
$ws.Cells.Item(121, 10).Value = 817.8
$ws.Cells.Item(121, 11).Value = 2510.4999
$ws.Cells.Item(121, 12).Value = 2453.4
$ws.Cells.Item(121, 13).Value = -1200.4999
$ws.Cells.Item(121, 14).Value = -5073.4
$ws.Cells.Item(129, 8).Value = 2910.8
$ws.Cells.Item(129, 9).Value = 1500
$ws.Cells.Item(129, 11).Value = 4500
$ws.Cells.Item(129, 13).Value = 500

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 3378.2727
$ws.Cells.Item(113, 9).Value = 2945.75
$ws.Cells.Item(113, 11).Value = 2945.75
$ws.Cells.Item(113, 13).Value = -775.75
$ws.Cells.Item(126, 8).Value = 2840.2307
$ws.Cells.Item(126, 9).Value = 2802.3333
$ws.Cells.Item(126, 11).Value = 8406.999899999999
$ws.Cells.Item(126, 13).Value = -5936.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 3111
$ws.Cells.Item(46, 9).Value = 2799.8
$ws.Cells.Item(46, 11).Value = 2799.8
$ws.Cells.Item(46, 13).Value = -2611.8
$ws.Cells.Item(132, 8).Value = 2465.125
$ws.Cells.Item(132, 9).Value = 2396.7058
$ws.Cells.Item(132, 11).Value = 7190.117400000001
$ws.Cells.Item(132, 13).Value = -4660.117400000001
$ws.Cells.Item(136, 8).Value = 5268.5557
$ws.Cells.Item(136, 9).Value = 6203.857
$ws.Cells.Item(136, 10).Value = 1995
$ws.Cells.Item(136, 11).Value = 18611.571
$ws.Cells.Item(136, 12).Value = 5985
$ws.Cells.Item(136, 13).Value = -16061.571
$ws.Cells.Item(136, 14).Value = -11085

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(45, 8).Value = 57174.8
$ws.Cells.Item(45, 9).Value = 17990
$ws.Cells.Item(45, 11).Value = 17990
$ws.Cells.Item(45, 13).Value = -17499
$ws.Cells.Item(121, 8).Value = 98999
$ws.Cells.Item(121, 10).Value = 98999
$ws.Cells.Item(121, 12).Value = 98999
$ws.Cells.Item(121, 14).Value = -102493
$ws.Cells.Item(132, 8).Value = 2031.2307
$ws.Cells.Item(132, 9).Value = 2310.1
$ws.Cells.Item(132, 10).Value = 1101.6666
$ws.Cells.Item(132, 11).Value = 6930.299999999999
$ws.Cells.Item(132, 12).Value = 3304.9998
$ws.Cells.Item(132, 13).Value = -4400.299999999999
$ws.Cells.Item(132, 14).Value = -8364.9998
